$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. SupIm sheet: extend the timeseries from 1 data row (t=1) to 12 data rows
#    (t=1..12), mirroring the existing row 3 pattern (A: incrementing index
#    with style "4", B/C: plain numbers, D: style "6").
# ---------------------------------------------------------------------------
$supim = $wb.Worksheets.Item("SupIm")

for ($row = 4; $row -le 14; $row++) {
    $t = $row - 2
    # Set the raw numeric values first (while the cells still use the
    # General format) so Excel stores them as numbers rather than text,
    # then copy the formatting (including number format/style) from row 3
    # afterwards so the stored value type is not coerced to text.
    $supim.Range("A$row").Value2 = $t
    $supim.Range("B$row").Value2 = 0.48099999999999998
    $supim.Range("C$row").Value2 = 0.3
    $supim.Range("D$row").Value2 = 0.20699999999999999
}

$supim.Range("A3:D3").Copy() | Out-Null
$supim.Range("A4:D14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Make SupIm the active sheet/tab, with K13 selected (was Process sheet
#    before, now switches to SupIm).
# ---------------------------------------------------------------------------
$supim.Activate() | Out-Null
$supim.Range("K13").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Process sheet: merge the two conditional formatting rules covering
#    A11:C13 into a single rule over the contiguous range A11:C13.
# ---------------------------------------------------------------------------
$process = $wb.Worksheets.Item("Process")

# Remove the separate rule that only applied to B11.
$process.Range("B11").FormatConditions.Delete() | Out-Null

# Extend the remaining rule (currently sqref="A12:C13 A11 C11") so it
# applies to the full contiguous block A11:C13, and make it the top
# (first) priority rule.
$fc = $process.Range("A11").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($process.Range("A11:C13"))
$fc = $process.Range("A11").FormatConditions.Item(1)
$fc.SetFirstPriority()

Write-Output "done"
